$wb = $excel.ActiveWorkbook

# The "AW & ST: App" reporting-type row (row 6) was dropped from the two
# per-sample summary exports. Deleting the row shifts the remaining
# reporting-type rows up by one and shrinks the used range by a row.

$ws1 = $wb.Worksheets.Item("Sampling_Delay_Mean_IQR")
$ws1.Rows("6").Delete()

$ws2 = $wb.Worksheets.Item("Cumulative_Sampling_Delay")
$ws2.Rows("6").Delete()
